$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "B2" = 0.8367346938775511; "C2" = 0.8367346938775511; "D2" = 0.8367346938775511; "E2" = 49
    "B3" = 0.8181818181818182; "C3" = 0.8181818181818182; "D3" = 0.8181818181818182; "E3" = 44
    "B4" = 0.8279569892473119; "C4" = 0.8279569892473119; "D4" = 0.8279569892473119; "E4" = 0.8279569892473119
    "B5" = 0.8274582560296846; "C5" = 0.8274582560296846; "D5" = 0.8274582560296846
    "B6" = 0.8279569892473119; "C6" = 0.8279569892473119; "D6" = 0.8279569892473119
    "B7" = 0.7777777777777778; "C7" = 0.7142857142857143; "D7" = 0.7446808510638298; "E7" = 49
    "B8" = 0.7083333333333334; "C8" = 0.7727272727272727; "D8" = 0.7391304347826088; "E8" = 44
    "B10" = 0.7430555555555556; "C10" = 0.7435064935064934; "D10" = 0.7419056429232193
    "B11" = 0.7449223416965354; "D11" = 0.7420548476619618
    "B12" = 0.8888888888888888; "C12" = 0.8163265306122449; "D12" = 0.851063829787234; "E12" = 49
    "B13" = 0.8125; "C13" = 0.8863636363636364; "D13" = 0.8478260869565218; "E13" = 44
    "B14" = 0.8494623655913979; "C14" = 0.8494623655913979; "D14" = 0.8494623655913979; "E14" = 0.8494623655913979
    "B15" = 0.8506944444444444; "C15" = 0.8513450834879406; "D15" = 0.8494449583718779
    "B16" = 0.8527479091995219; "C16" = 0.8494623655913979; "D16" = 0.8495319944694777
    "B17" = 0.8222222222222222; "C17" = 0.7551020408163265; "D17" = 0.7872340425531914; "E17" = 49
    "B18" = 0.75; "C18" = 0.8181818181818182; "D18" = 0.7826086956521738; "E18" = 44
    "B19" = 0.7849462365591398; "C19" = 0.7849462365591398; "D19" = 0.7849462365591398; "E19" = 0.7849462365591398
    "B20" = 0.7861111111111111; "C20" = 0.7866419294990723; "D20" = 0.7849213691026826
    "B21" = 0.78805256869773; "C21" = 0.7849462365591398; "D21" = 0.7850457063849681
    "B22" = 0.8823529411764706; "C22" = 0.9183673469387755; "D22" = 0.9; "E22" = 49
    "B23" = 0.9047619047619048; "C23" = 0.8636363636363636; "D23" = 0.8837209302325582; "E23" = 44
    "B24" = 0.8924731182795699; "C24" = 0.8924731182795699; "D24" = 0.8924731182795699; "E24" = 0.8924731182795699
    "B25" = 0.8935574229691876; "C25" = 0.8910018552875696; "D25" = 0.8918604651162791
    "B26" = 0.8929550314749556; "C26" = 0.8924731182795699; "D26" = 0.8922980745186295
}

foreach ($key in $values.Keys) {
    $ws.Range($key).Value = $values[$key]
}
